$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: round coordinates, drop empty time columns ---
$ws.Range("Q54").Value = 623280
$ws.Range("R54").Value = 6951831
$ws.Range("Z54").ClearContents()
$ws.Range("AB54").ClearContents()

# --- Row 55: becomes the former row 56 content (rows 55/56 swapped) ---
$ws.Range("A55").Value = 112128672
$ws.Range("B55").Value = 90689
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 5966
$ws.Range("F55").Value = "Motaggsvamp"
$ws.Range("G55").Value = "Sarcodon squamosus"
$ws.Range("H55").Value = "(Schaeff.) Quél."
$ws.Range("I55").ClearContents()
$ws.Range("J55").ClearContents()
$ws.Range("Q55").Value = 623280
$ws.Range("R55").Value = 6951831
$ws.Range("Z55").ClearContents()
$ws.Range("AB55").ClearContents()
$ws.Range("AC55").Value = "stgen upp mot vägen"

# --- Row 56: becomes the former row 55 content ---
$ws.Range("A56").Value = 112128712
$ws.Range("B56").Value = 88914
$ws.Range("D56").Value = "VU"
$ws.Range("E56").Value = 2051
$ws.Range("F56").Value = "Rotfingersvamp"
$ws.Range("G56").Value = "Ramaria boreimaxima"
$ws.Range("H56").Value = "Kytöv. & M.Toivonen"
$ws.Range("I56").NumberFormat = "@"
$ws.Range("I56").Value = "3"
$ws.Range("I56").Style = "Normal"
$ws.Range("J56").Value = "fruktkroppar"
$ws.Range("Q56").Value = 623280
$ws.Range("R56").Value = 6951831
$ws.Range("Z56").ClearContents()
$ws.Range("AB56").ClearContents()
$ws.Range("AC56").Value = "stigen upp mot vägen"
